$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-looking numeric price cells to stay as Text (matches original inlineStr data type)
$textCells = @("D5","D8","D11","D16","D18","D22","D24","D25","D35","D36","D39","D40","D43","D48","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range('D2').Value = '27.961.99'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '1.634.42'
$ws.Range('E3').Value = '  -0.50%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '211.99'
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '23.45'
$ws.Range('E8').Value = '  -2.02%  '
$ws.Range('E9').Value = '  -2.15%  '
$ws.Range('D11').Value = '0.0880'
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').Value = '1.865.62'
$ws.Range('E12').Value = '  -0.55%  '
$ws.Range('D13').Value = '1.628.04'
$ws.Range('E13').Value = '  -0.88%  '
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('E15').Value = '  -2.26%  '
$ws.Range('D16').Value = '65.79'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').Value = '27.960.42'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').Value = '231.77'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').Value = '10.42'
$ws.Range('E22').Value = '  -5.94%  '
$ws.Range('E23').Value = '  -0.79%  '
$ws.Range('D24').Value = '2.07'
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('D25').Value = '155.25'
$ws.Range('E25').Value = '  +1.85%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('E27').Value = '  -0.66%  '
$ws.Range('E28').Value = '  -0.81%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('E30').Value = '  -0.50%  '
$ws.Range('E31').Value = '  -1.02%  '
$ws.Range('E32').Value = '  +1.65%  '
$ws.Range('E33').Value = '  -0.86%  '
$ws.Range('D34').Value = '1.405.87'
$ws.Range('E34').Value = '  -1.25%  '
$ws.Range('B35').Value = 'TrustWalletToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D35').Value = '1.05'
$ws.Range('E35').Value = '  +13.45%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').Value = '1.57'
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('E38').Value = '  +2.10%  '
$ws.Range('D39').Value = '0.555'
$ws.Range('E39').Value = '  -0.65%  '
$ws.Range('D40').Value = '0.867'
$ws.Range('E40').Value = '  -2.85%  '
$ws.Range('E41').Value = '  -0.71%  '
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').Value = '66.92'
$ws.Range('E43').Value = '  -0.58%  '
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('E46').Value = '  -0.33%  '
$ws.Range('D47').Value = '1.776.26'
$ws.Range('E47').Value = '  -0.45%  '
$ws.Range('D48').Value = '88.28'
$ws.Range('E48').Value = '  -0.73%  '
$ws.Range('E49').Value = '  -0.62%  '
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('D51').Value = '0.0505'
$ws.Range('E51').Value = '  -0.30%  '
